$d = $word.ActiveDocument
$sec = $d.Sections.First

# Footer 1 (primary footer) and Footer 2 (first-page footer) each carry the
# Pearson Edexcel logo picture. Its shape is currently named "image2.png" in
# both places; rename it to "image1.png".
for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        foreach ($s in $ftr.Range.InlineShapes) {
            if ($s.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $s.Name = "image1.png"
            }
        }
    }
}

# Header 1 (first-page header) carries the BTec logo picture. Its shape is
# currently named "image1.jpg"; rename it to "image2.jpg".
for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        foreach ($s in $hdr.Range.InlineShapes) {
            if ($s.AlternativeText -eq "BTec_Logo-Orange") {
                $s.Name = "image2.jpg"
            }
        }
    }
}
